$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D5").Value = "'19000"
$ws.Range("D6").Value = "'6000"
